# Update: Threat Alert Report - 2026-01-25 09:07
#
# Refreshes the TUU-ATZ threat alert sheet: rewrites the date/fare/baggage
# figures for existing rows, retires the "MEDIUM THREAT - MONITOR" banding
# (style+fill) since nothing uses it anymore, appends two new observations
# (rows 7-8), and narrows the IMPACT column now that it only ever shows
# "LOW THREAT".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a date-like label into column A without Excel silently
# re-interpreting "27-FEB-26" as a date serial. Building it as a text
# formula and then pasting-by-value keeps the stored cell a plain string.
# ---------------------------------------------------------------------
function Set-DateText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Formula = '="' + $text + '"'
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Apply the formatting (fill/border/font/alignment) of a known-good cell
# onto a target cell without disturbing the value already written there.
function Copy-Format($srcRef, $dstRef) {
    $ws.Range($srcRef).Copy() | Out-Null
    $ws.Range($dstRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

function Set-DataRow($r, $date, $flight, $airline, $oalFare, $ourFare, $fareDif, $oalBag, $ourBag, $bagDif, $impact, $currency) {
    Set-DateText "A$r" $date
    $ws.Range("B$r").Value = $flight
    $ws.Range("C$r").Value = $airline
    $ws.Range("D$r").Value = $oalFare
    $ws.Range("E$r").Value = $ourFare
    $ws.Range("F$r").Value = $fareDif
    $ws.Range("G$r").Value = $oalBag
    $ws.Range("H$r").Value = $ourBag
    $ws.Range("I$r").Value = $bagDif
    $ws.Range("J$r").Value = $impact
    $ws.Range("K$r").Value = $currency

    # Re-stamp the row's look (fill/border/font/alignment) from the
    # template data row, since new rows start with no style at all and
    # Set-DateText's value-only paste above strips A's style too.
    Copy-Format "A2:I2" "A${r}:I${r}"
    Copy-Format "J2" "J$r"
    Copy-Format "K2" "K$r"
}

# Row 2 — figures refreshed, wording/threat level unchanged
Set-DataRow 2 "30-JAN-26" "SM-328" "EgyptAir MS-812" 431 427 4 46 30 -16 "LOW THREAT" "SAR"

# Row 3 — now EgyptAir (was Air Arabia Egypt), new date/fares
Set-DataRow 3 "06-FEB-26" "SM-328" "EgyptAir MS-812" 383 517 -134 46 30 -16 "LOW THREAT" "SAR"

# Row 4
Set-DataRow 4 "13-FEB-26" "SM-328" "EgyptAir MS-812" 383 517 -134 46 30 -16 "LOW THREAT" "SAR"

# Row 5 — threat level downgraded from MEDIUM to LOW
Set-DataRow 5 "20-FEB-26" "SM-328" "EgyptAir MS-812" 383 406 -23 46 30 -16 "LOW THREAT" "SAR"

# Row 6 — now Air Arabia Egypt (was EgyptAir)
Set-DataRow 6 "27-FEB-26" "SM-328" "Air Arabia Egypt E5-590" 399 602 -203 30 30 0 "LOW THREAT" "SAR"

# Rows 7-8 are brand new observations appended to the report
Set-DataRow 7 "27-FEB-26" "SM-328" "EgyptAir MS-812" 516 602 -86 46 30 -16 "LOW THREAT" "SAR"
Set-DataRow 8 "13-MAR-26" "SM-328" "EgyptAir MS-812" 917 932 -15 46 30 -16 "LOW THREAT" "SAR"

# The IMPACT column no longer needs to be wide enough for
# "MEDIUM THREAT - MONITOR" now that every row reads "LOW THREAT".
$ws.Columns.Item(10).ColumnWidth = 11.17

Write-Host "sheet refreshed"
